# Scheduled market-data refresh: update currentAveragePrice / Leve profit
# columns (H:N) per leve row across the crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 6720.7295
$ws.Range("I76").Value = 7427.24
$ws.Range("J76").Value = 5248.8335
$ws.Range("K76").Value = 7427.24
$ws.Range("L76").Value = 5248.8335
$ws.Range("M76").Value = -7112.24
$ws.Range("N76").Value = -5878.8335
# Row 79
$ws.Range("H79").Value = 6720.7295
$ws.Range("I79").Value = 7427.24
$ws.Range("J79").Value = 5248.8335
$ws.Range("K79").Value = 7427.24
$ws.Range("L79").Value = 5248.8335
$ws.Range("M79").Value = -6335.24
$ws.Range("N79").Value = -7432.8335
# Row 100
$ws.Range("H100").Value = 5304.72
$ws.Range("I100").Value = 2435.5
$ws.Range("J100").Value = 6654.9414
$ws.Range("K100").Value = 2435.5
$ws.Range("L100").Value = 6654.9414
$ws.Range("M100").Value = -1894.5
$ws.Range("N100").Value = -7736.9414
# Row 135
$ws.Range("H135").Value = 1060.0454
$ws.Range("I135").Value = 553.7059
$ws.Range("K135").Value = 4983.3531
$ws.Range("M135").Value = -2448.3531
# Row 137
$ws.Range("H137").Value = 30257.723
$ws.Range("I137").Value = 42786.64
$ws.Range("J137").Value = 1782.909
$ws.Range("K137").Value = 128359.92
$ws.Range("L137").Value = 5348.727000000001
$ws.Range("M137").Value = -125809.92
$ws.Range("N137").Value = -10448.727
$ws = $wb.Worksheets.Item("ARM")
# Row 12
$ws.Range("H12").Value = 1000
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1346
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4809.125
$ws.Range("I86").Value = 4773.375
$ws.Range("J86").Value = 4844.875
$ws.Range("K86").Value = 4773.375
$ws.Range("L86").Value = 4844.875
$ws.Range("M86").Value = -3650.375
$ws.Range("N86").Value = -7090.875
# Row 89
$ws.Range("H89").Value = 4809.125
$ws.Range("I89").Value = 4773.375
$ws.Range("J89").Value = 4844.875
$ws.Range("K89").Value = 23866.875
$ws.Range("L89").Value = 24224.375
$ws.Range("M89").Value = -18250.875
$ws.Range("N89").Value = -35456.375
# Row 134
$ws.Range("H134").Value = 43909.69
$ws.Range("I134").Value = 49166
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 147498
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -144963
$ws.Range("N134").Value = -50070
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3829.611
$ws.Range("I31").Value = 3605.5386
$ws.Range("J31").Value = 3956.261
$ws.Range("K31").Value = 3605.5386
$ws.Range("L31").Value = 3956.261
$ws.Range("M31").Value = -3310.5386
$ws.Range("N31").Value = -4546.261
# Row 34
$ws.Range("H34").Value = 3829.611
$ws.Range("I34").Value = 3605.5386
$ws.Range("J34").Value = 3956.261
$ws.Range("K34").Value = 3605.5386
$ws.Range("L34").Value = 3956.261
$ws.Range("M34").Value = -3403.5386
$ws.Range("N34").Value = -4360.261
# Row 58
$ws.Range("H58").Value = 1683.1305
$ws.Range("I58").Value = 832.6
$ws.Range("J58").Value = 2337.3845
$ws.Range("K58").Value = 832.6
$ws.Range("L58").Value = 2337.3845
$ws.Range("M58").Value = -629.6
$ws.Range("N58").Value = -2743.3845
# Row 134
$ws.Range("H134").Value = 1598.2609
$ws.Range("I134").Value = 924.7586
$ws.Range("J134").Value = 2747.1765
$ws.Range("K134").Value = 2774.2758
$ws.Range("L134").Value = 8241.529500000001
$ws.Range("M134").Value = -239.2757999999999
$ws.Range("N134").Value = -13311.5295
# Row 136
$ws.Range("H136").Value = 1683.1305
$ws.Range("I136").Value = 832.6
$ws.Range("J136").Value = 2337.3845
$ws.Range("K136").Value = 2497.8
$ws.Range("L136").Value = 7012.1535
$ws.Range("M136").Value = 52.19999999999982
$ws.Range("N136").Value = -12112.1535
$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 880.9524
$ws.Range("I12").Value = 880.9524
$ws.Range("K12").Value = 880.9524
$ws.Range("M12").Value = -740.9524
# Row 113
$ws.Range("H113").Value = 5328.478
$ws.Range("I113").Value = 5827.8
$ws.Range("J113").Value = 1999.6666
$ws.Range("K113").Value = 5827.8
$ws.Range("L113").Value = 1999.6666
$ws.Range("M113").Value = -3657.8
$ws.Range("N113").Value = -6339.6666
# Row 126
$ws.Range("H126").Value = 2848.8
$ws.Range("I126").Value = 2679.9333
$ws.Range("J126").Value = 3355.4
$ws.Range("K126").Value = 8039.7999
$ws.Range("L126").Value = 10066.2
$ws.Range("M126").Value = -5569.7999
$ws.Range("N126").Value = -15006.2
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2404.9443
$ws.Range("I7").Value = 2209
$ws.Range("J7").Value = 2712.8572
$ws.Range("K7").Value = 2209
$ws.Range("L7").Value = 2712.8572
$ws.Range("M7").Value = -2097
$ws.Range("N7").Value = -2936.8572
# Row 16
$ws.Range("H16").Value = 1508.5652
$ws.Range("J16").Value = 1013.7143
$ws.Range("L16").Value = 1013.7143
$ws.Range("N16").Value = -1353.7143
# Row 93
$ws.Range("H93").Value = 2714.2856
$ws.Range("I93").Value = 2600
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2600
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1352
$ws.Range("N93").Value = -5496
# Row 126
$ws.Range("H126").Value = 2404.9443
$ws.Range("I126").Value = 2209
$ws.Range("J126").Value = 2712.8572
$ws.Range("K126").Value = 6627
$ws.Range("L126").Value = 8138.571599999999
$ws.Range("M126").Value = -4157
$ws.Range("N126").Value = -13078.5716
# Row 132
$ws.Range("H132").Value = 10235.615
$ws.Range("I132").Value = 2845.2666
$ws.Range("J132").Value = 20313.363
$ws.Range("K132").Value = 8535.799800000001
$ws.Range("L132").Value = 60940.08900000001
$ws.Range("M132").Value = -6005.799800000001
$ws.Range("N132").Value = -66000.08900000001
# Row 133
$ws.Range("H133").Value = 38990
$ws.Range("J133").Value = 38990
$ws.Range("L133").Value = 38990
$ws.Range("N133").Value = -44050
# Row 136
$ws.Range("H136").Value = 7100.375
$ws.Range("I136").Value = 4293.2144
$ws.Range("J136").Value = 11030.4
$ws.Range("K136").Value = 12879.6432
$ws.Range("L136").Value = 33091.2
$ws.Range("M136").Value = -10329.6432
$ws.Range("N136").Value = -38191.2
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 24342.512
$ws.Range("I132").Value = 39844.69
$ws.Range("J132").Value = 3129
$ws.Range("K132").Value = 119534.07
$ws.Range("L132").Value = 9387
$ws.Range("M132").Value = -117004.07
$ws.Range("N132").Value = -14447
